# Results from July 28, 2020 10:32:58 PM America/Los_Angeles TZ run
# Updates covid_disparities_output sheet with the latest scrape results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Texas -- Bexar County ---
$ws.Cells.Item(3, 2).Value = 44040
$ws.Cells.Item(3, 3).Value = 37984
$ws.Cells.Item(3, 4).Value = 335

# --- Row 6: California - San Francisco (error message only) ---
$ws.Cells.Item(6, 15).Value = 'An error occurred. ... KeyError(''Date_Uploaded.Data as of'')'

# --- Row 7: Tennessee ---
$ws.Cells.Item(7, 2).Value = 44040
$ws.Cells.Item(7, 3).Value = 99044
$ws.Cells.Item(7, 4).Value = 999
$ws.Cells.Item(7, 5).Value = 18749
$ws.Cells.Item(7, 6).Value = 343
$ws.Cells.Item(7, 7).Value = 18.93
$ws.Cells.Item(7, 8).Value = 34.33

# --- Row 11: California - San Diego ---
$ws.Cells.Item(11, 2).Value = 44040
$ws.Cells.Item(11, 3).Value = 28005
$ws.Cells.Item(11, 4).Value = 547
$ws.Cells.Item(11, 5).Value = 1034
$ws.Cells.Item(11, 7).Value = 4.73
$ws.Cells.Item(11, 8).Value = 3.77
$ws.Cells.Item(11, 11).Value = 21845
$ws.Cells.Item(11, 12).Value = 531

# --- Row 36: Washington ---
$ws.Cells.Item(36, 2).Value = 44040
$ws.Cells.Item(36, 3).Value = 54205
$ws.Cells.Item(36, 4).Value = 1548
$ws.Cells.Item(36, 5).Value = 1984
$ws.Cells.Item(36, 7).Value = 5.52
$ws.Cells.Item(36, 8).Value = 3.26
$ws.Cells.Item(36, 11).Value = 35958
$ws.Cells.Item(36, 12).Value = 1505

# --- Row 39: Delaware (error message only) ---
$ws.Cells.Item(39, 15).Value = 'An error occurred. ... WebDriverException(''unknown error: session deleted because of page crash\nfrom unknown error: cannot determine loading status\nfrom tab crashed\n  (Session info: headless chrome=83.0.4103.116)'', None, None)'

# --- Row 41: Iowa ---
$ws.Cells.Item(41, 2).Value = 44041
$ws.Cells.Item(41, 3).Value = 42928
$ws.Cells.Item(41, 5).Value = 3407
$ws.Cells.Item(41, 7).Value = 7.94

# --- Row 45: Ohio -- scrape failed this run, clear the scraped columns ---
$ws.Range("B45:H45").ClearContents()
$ws.Range("B45:H45").Style = "Normal"
$ws.Cells.Item(45, 10).Value = $false
$ws.Range("K45:L45").ClearContents()
$ws.Range("K45:L45").Style = "Normal"
$ws.Cells.Item(45, 15).Value = 'An error occurred. ... AttributeError("''NoneType'' object has no attribute ''body''")'

# --- Row 50: NewYork -- scrape failed this run, clear the scraped columns ---
$ws.Range("B50:D50").ClearContents()
$ws.Range("B50:D50").Style = "Normal"
$ws.Cells.Item(50, 6).ClearContents()
$ws.Cells.Item(50, 6).Style = "Normal"
$ws.Cells.Item(50, 8).ClearContents()
$ws.Cells.Item(50, 8).Style = "Normal"
$ws.Cells.Item(50, 12).ClearContents()
$ws.Cells.Item(50, 12).Style = "Normal"
$ws.Cells.Item(50, 15).Value = 'An error occurred. ... ConnectionRefusedError(111, ''Connection refused'')'

# --- Insert two new rows (51, 52) before the old SouthCarolina row, which ---
# --- shifts SouthCarolina from row 51 down to row 53.                    ---
$ws.Range("A51:A52").EntireRow.Insert()

# --- Row 51 (new): Wyoming ---
$ws.Cells.Item(51, 1).Value = "Wyoming"
$ws.Range("B51:H51").ClearContents()
$ws.Range("B51:H51").Style = "Normal"
$ws.Cells.Item(51, 9).Value = $false
$ws.Cells.Item(51, 10).Value = $false
$ws.Range("K51:L51").ClearContents()
$ws.Range("K51:L51").Style = "Normal"
$ws.Cells.Item(51, 13).Value = 5540
$ws.Cells.Item(51, 14).Value = 0.95
$ws.Cells.Item(51, 15).Value = 'An error occurred. ... JSONDecodeError(''Expecting value: line 1 column 1 (char 0)'')'

# --- Row 52 (new): SouthDakota ---
$ws.Cells.Item(52, 1).Value = "SouthDakota"
$ws.Cells.Item(52, 2).Value = 44040
$ws.Cells.Item(52, 3).Value = 8492
$ws.Cells.Item(52, 4).Value = 123
$ws.Cells.Item(52, 5).Value = 1008
$ws.Cells.Item(52, 6).ClearContents()
$ws.Cells.Item(52, 6).Style = "Normal"
$ws.Cells.Item(52, 7).Value = 11.87
$ws.Cells.Item(52, 8).ClearContents()
$ws.Cells.Item(52, 8).Style = "Normal"
$ws.Cells.Item(52, 9).Value = $false
$ws.Cells.Item(52, 10).Value = $false
$ws.Cells.Item(52, 11).Value = 8492
$ws.Cells.Item(52, 12).ClearContents()
$ws.Cells.Item(52, 12).Style = "Normal"
$ws.Cells.Item(52, 13).ClearContents()
$ws.Cells.Item(52, 13).Style = "Normal"
$ws.Cells.Item(52, 14).ClearContents()
$ws.Cells.Item(52, 14).Style = "Normal"
$ws.Cells.Item(52, 15).Value = "Success!"

Write-Output "done"
